$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting from B1 into the new header cells C1:D1, then set values
$ws.Range("B1").Copy()
$ws.Range("C1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C1").Value = "Postal"
$ws.Range("D1").Value = "Population"

# Postal abbreviations (column C) and population figures (column D)
$postal = @("AL","AK","AZ","AR","CA","CO","CT","DE","DC","FL","GA","HI","ID","IL","IN","IA","KS","KY","LA","ME","MD","MA","MI","MN","MS","MO","MT","NE","NV","NH","NJ","NM","NY","NC","ND","OH","OK","OR","PA","PR","RI","SC","SD","TN","TX","UT","VT","VA","WA","WV","WI","WY")
$population = @(4849377,736732,6731484,2966369,38802500,5355866,3596677,935614,658893,19893297,10097343,1419561,1634464,12880580,6596855,3107126,2904021,4413457,4649676,1330089,5976407,6745408,9909877,5457173,2994079,6063589,1023579,1881503,2839098,1326813,8938175,2085572,19746227,9943964,739482,11594163,3878051,3970239,12787209,3548397,1055173,4832482,853175,6549352,26956958,2942902,626562,8326289,7061530,1850326,5757564,584153)

for ($i = 0; $i -lt $postal.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $postal[$i]
    $ws.Cells.Item($row, 4).Value = $population[$i]
}
